$d = $word.ActiveDocument
$t = $d.Tables(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $target = $d.Range($r.Start, $r.End - 1)
    $target.Text = $newText
}

function Set-CellFontSize($table, $row, $col, $points) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $target = $d.Range($r.Start, $r.End - 1)
    $target.Font.Size = $points
}

function Set-CellAlignRight($table, $row, $col) {
    $cell = $table.Cell($row, $col)
    $para = $cell.Range.Paragraphs(1)
    $para.Alignment = 2
}

# 1) Header row font size 11pt -> 10pt (sz 22 -> 20)
Set-CellFontSize $t 1 1 10
Set-CellFontSize $t 1 2 10
Set-CellFontSize $t 1 3 10
Set-CellFontSize $t 1 4 10
Set-CellFontSize $t 1 5 10
Set-CellFontSize $t 1 6 10
Set-CellFontSize $t 1 7 10
Set-CellFontSize $t 1 8 10

# 2) STOCK-name column paragraph alignment left -> right (rows 2..12)
Set-CellAlignRight $t 2 1
Set-CellAlignRight $t 3 1
Set-CellAlignRight $t 4 1
Set-CellAlignRight $t 5 1
Set-CellAlignRight $t 6 1
Set-CellAlignRight $t 7 1
Set-CellAlignRight $t 8 1
Set-CellAlignRight $t 9 1
Set-CellAlignRight $t 10 1
Set-CellAlignRight $t 11 1
Set-CellAlignRight $t 12 1

# 3) Numeric cell text reformatting
Set-CellText $t 2 2 "-5457.000"
Set-CellText $t 2 3 "-5495.000"
Set-CellText $t 2 4 "0.002"
Set-CellText $t 2 5 "-2505.000"
Set-CellText $t 2 6 "154496.552"
Set-CellText $t 3 2 "-1207.300"
Set-CellText $t 3 3 "-1207.000"
Set-CellText $t 3 4 "0.008"
Set-CellText $t 3 5 "-2055.500"
Set-CellText $t 3 6 "19008.038"
Set-CellText $t 3 8 "60174.885"
Set-CellText $t 4 2 "-3956.000"
Set-CellText $t 4 3 "-3956.000"
Set-CellText $t 4 5 "-2518.800"
Set-CellText $t 4 6 "117668.826"
Set-CellText $t 4 8 "34368.826"
Set-CellText $t 5 2 "15249.000"
Set-CellText $t 5 3 "-1086.000"
Set-CellText $t 5 4 "0.016"
Set-CellText $t 5 5 "-4731.000"
Set-CellText $t 5 6 "393224.787"
Set-CellText $t 6 2 "-1396.800"
Set-CellText $t 6 3 "-797.300"
Set-CellText $t 6 4 "0.000"
Set-CellText $t 6 5 "-156.600"
Set-CellText $t 6 6 "61655.065"
Set-CellText $t 6 8 "16189.087"
Set-CellText $t 7 2 "-1765.000"
Set-CellText $t 7 4 "0.010"
Set-CellText $t 7 5 "-1675.000"
Set-CellText $t 7 6 "387884.615"
Set-CellText $t 8 2 "-984.100"
Set-CellText $t 8 3 "-728.100"
Set-CellText $t 8 5 "-4386.300"
Set-CellText $t 8 6 "272582.439"
Set-CellText $t 9 2 "-498.240"
Set-CellText $t 9 3 "-954.600"
Set-CellText $t 9 4 "0.042"
Set-CellText $t 9 5 "-646.000"
Set-CellText $t 9 6 "15694.472"
Set-CellText $t 9 8 "13069.083"
Set-CellText $t 10 2 "-568.000"
Set-CellText $t 10 3 "-363.740"
Set-CellText $t 10 4 "0.460"
Set-CellText $t 10 5 "-209.590"
Set-CellText $t 10 6 "10048.103"
Set-CellText $t 10 8 "6070.205"
Set-CellText $t 11 4 "0.003"
Set-CellText $t 11 6 "17766.997"
Set-CellText $t 11 8 "13783.752"
Set-CellText $t 12 2 "-4901.900"
Set-CellText $t 12 3 "-4915.800"
Set-CellText $t 12 4 "0.010"
Set-CellText $t 12 5 "-11746.000"
Set-CellText $t 12 6 "128653.200"
